$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Row 17 holds the "repaymentstrategy" input; update its value from
# "RBI (India)" to "Overdue/Due Fee/Int,Principal"
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Reflect the active cell selection as it was left after the edit
$ws.Activate()
$ws.Range("B17").Select()
